$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new feed-log rows (150 and 151) after the existing data (ends at row 149)
$ws.Range("A150").Value = 149
$ws.Range("B150").Value = 1
$ws.Range("C150").Value = "2024-06-18 03:14:56"
$ws.Range("D150").Value = 200
$ws.Range("E150").Value = 13

$ws.Range("A151").Value = 150
$ws.Range("B151").Value = 2
$ws.Range("C151").Value = "2024-06-18 03:14:56"
$ws.Range("D151").Value = 200
$ws.Range("E151").Value = 0
